$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-51: (row, Coin(B), Link(C), Price(D), Volume(E), DNeedsTextPrefix)
#
# Transcribed from the target OOXML diff. Row 34 gains a newly inserted coin
# ("Frax"), which pushes the coin that used to occupy each subsequent row down
# by one; the coin that used to be the last row (Decentraland) falls off the
# bottom of the A1:E51 range and is dropped. "EnergySwap" (row 49) keeps its
# rank but updates its price/volume, per the diff.
#
# The DNeedsTextPrefix flag marks Price values that look like plain numbers to
# Excel (e.g. "1.000", "247.06") so we can force Excel to store them as text
# (matching the source file's inlineStr cells) instead of silently normalizing
# them into numeric values and losing formatting such as trailing zeros.
$data = @(
    @(2, 'Bitcoin', 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc', '26.499.11', '  -0.31%  ', 0),
    @(3, 'Ethereum', 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth', '1.731.38', '  -0.68%  ', 0),
    @(4, 'TetherUSD', 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt', '1.000', '  +0.03%  ', 1),
    @(5, 'BNB', 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb', '247.06', '  +0.23%  ', 1),
    @(6, 'USDC', 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc', '1.000', '  -0.01%  ', 1),
    @(7, 'XRP', 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp', '0.4876', '  +1.21%  ', 1),
    @(8, 'Cardano', 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada', '0.2669', '  -0.91%  ', 1),
    @(9, 'Dogecoin', 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge', '0.06218', '  -0.70%  ', 1),
    @(10, 'WrappedEther', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', '1.732.17', '  -0.63%  ', 0),
    @(11, 'TRON', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx', '0.07068', '  -0.90%  ', 1),
    @(12, 'Solana', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol', '15.69', '  -0.88%  ', 1),
    @(13, 'Polkadot', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', '4.636', '  +2.68%  ', 1),
    @(14, 'Polygon', 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic', '0.6090', '  -2.12%  ', 1),
    @(15, 'Litecoin', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', '77.35', '  -0.04%  ', 1),
    @(16, 'Dai', 'https://coinranking.com/coin/MoTuySvg7+dai-dai', '1.000', '  -0.01%  ', 1),
    @(17, 'WrappedBTC', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', '26.498.55', '  -0.31%  ', 0),
    @(18, 'BinanceUSD', 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd', '1.000', '  -0.04%  ', 1),
    @(19, 'ShibaInu', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib', '0.000007226', '  +4.65%  ', 1),
    @(20, 'Avalanche', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', '11.52', '  -2.14%  ', 1),
    @(21, 'WrappedliquidstakedEther2.0', 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth', '1.955.94', '  -0.64%  ', 0),
    @(22, 'Uniswap', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', '4.524', '  -1.99%  ', 1),
    @(23, 'Cosmos', 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom', '8.764', '  -1.24%  ', 1),
    @(24, 'Chainlink', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', '5.250', '  -2.19%  ', 1),
    @(25, 'Monero', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', '139.00', '  +2.24%  ', 1),
    @(26, 'EthereumClassic', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', '15.43', '  +0.47%  ', 1),
    @(27, 'LidoDAOToken', 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo', '1.776', '  -1.99%  ', 1),
    @(28, 'Toncoin', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', '1.406', '  -1.94%  ', 1),
    @(29, 'BitcoinCash', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', '108.00', '  +1.18%  ', 1),
    @(30, 'InternetComputer(DFINITY)', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', '3.971', '  -0.99%  ', 1),
    @(31, 'Stellar', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', '0.08038', '  +1.87%  ', 1),
    @(32, 'Filecoin', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', '3.686', '  -1.41%  ', 1),
    @(33, 'Hedera', 'https://coinranking.com/coin/jad286TjB+hedera-hbar', '0.04570', '  -0.15%  ', 1),
    @(34, 'Frax', 'https://coinranking.com/coin/KfWtaeV1W+frax-frax', '1.000', '  +0.03%  ', 1),
    @(35, 'HuobiToken', 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht', '2.615', '  -0.09%  ', 1),
    @(36, 'ARBITRUM', 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb', '1.007', '  +0.79%  ', 1),
    @(37, 'ImmutableX', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', '0.6384', '  +0.14%  ', 1),
    @(38, 'TrustWalletToken', 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt', '0.8999', '  -3.57%  ', 1),
    @(39, 'RenderToken', 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr', '2.027', '  +2.06%  ', 1),
    @(40, 'MXToken', 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx', '2.399', '  -1.59%  ', 1),
    @(41, 'PaxDollar', 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp', '1.004', '  +0.14%  ', 1),
    @(42, 'VeChain', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet', '0.01506', '  -0.34%  ', 1),
    @(43, 'Quant', 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt', '101.35', '  -10.77%  ', 1),
    @(44, 'FraxShare', 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs', '5.447', '  -5.81%  ', 1),
    @(45, 'TheSandbox', 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand', '0.3894', '  -0.46%  ', 1),
    @(46, 'Aptos', 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt', '6.967', '  +2.99%  ', 1),
    @(47, 'Algorand', 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo', '0.1185', '  -2.30%  ', 1),
    @(48, 'Cronos', 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro', '0.05386', '  +0.97%  ', 1),
    @(49, 'EnergySwap', 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens', '7.825', '  -1.37%  ', 1),
    @(50, 'Elrond', 'https://coinranking.com/coin/omwkOTglq+elrond-egld', '30.54', '  -0.72%  ', 1),
    @(51, 'NEARProtocol', 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near', '1.247', '  -1.59%  ', 1)

)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    if ($row[5] -eq 1) {
        $ws.Cells.Item($r, 4).Value = "'" + $row[3]
    } else {
        $ws.Cells.Item($r, 4).Value = $row[3]
    }
    $ws.Cells.Item($r, 5).Value = $row[4]
}
